$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 6444
$ws.Range("C2").Value = 772
$ws.Range("D2").Value = 1490
$ws.Range("E2").Value = 75.7088764742396
$ws.Range("F2").Value = 78.78238341968913
$ws.Range("G2").Value = 86.29060402684564
$ws.Range("H2").Value = 4254
$ws.Range("I2").Value = 0.8719571687423648
$ws.Range("J2").Value = 3034
$ws.Range("K2").Value = 4.988490628082867
$ws.Range("L2").Value = 3383
$ws.Range("M2").Value = 2.63119006323256
$ws.Range("N2").Value = 12
$ws.Range("O2").Value = 0.002459681717185796
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 0.003288391976323578
$ws.Range("R2").Value = 3
$ws.Range("S2").Value = 0.00233330481516337

# Row 6
$ws.Range("B6").Value = 5936
$ws.Range("C6").Value = 383
$ws.Range("D6").Value = 767
$ws.Range("E6").Value = 105.8586590296496
$ws.Range("F6").Value = 66.59007832898172
$ws.Range("G6").Value = 66.94263363754889
$ws.Range("H6").Value = 1668
$ws.Range("I6").Value = 0.2654457435584052
$ws.Range("J6").Value = 1644
$ws.Range("K6").Value = 6.446047678795483
$ws.Range("L6").Value = 1940
$ws.Range("M6").Value = 3.778362060570649
$ws.Range("N6").Value = 227
$ws.Range("O6").Value = 0.03612481042431534
$ws.Range("P6").Value = 10
$ws.Range("Q6").Value = 0.03920953575909662
$ws.Range("R6").Value = 9
$ws.Range("S6").Value = 0.0175284837861525

# Row 11
$ws.Range("B11").Value = 7926
$ws.Range("C11").Value = 1132
$ws.Range("D11").Value = 2266
$ws.Range("E11").Value = 47.81175876860964
$ws.Range("F11").Value = 47.93109540636043
$ws.Range("G11").Value = 47.40644307149162
$ws.Range("H11").Value = 2738
$ws.Range("I11").Value = 0.7225113205754757
$ws.Range("J11").Value = 1861
$ws.Range("K11").Value = 3.429908953518375
$ws.Range("L11").Value = 2215
$ws.Range("M11").Value = 2.061942042206976
$ws.Range("N11").Value = 6
$ws.Range("O11").Value = 0.001583297269339976
$ws.Range("P11").Value = 1
$ws.Range("Q11").Value = 0.001843046186737439
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0

# Row 12
$ws.Range("B12").Value = 4432
$ws.Range("C12").Value = 633
$ws.Range("D12").Value = 1267
$ws.Range("E12").Value = 23.67170577617329
$ws.Range("F12").Value = 21.32859399684044
$ws.Range("G12").Value = 22.94869771112865
$ws.Range("H12").Value = 3422
$ws.Range("I12").Value = 3.261750212080486
$ws.Range("J12").Value = 2088
$ws.Range("K12").Value = 15.46552107251315
$ws.Range("L12").Value = 2591
$ws.Range("M12").Value = 8.911129453845096
$ws.Range("N12").Value = 3
$ws.Range("O12").Value = 0.002859512167224272

# Row 14
$ws.Range("B14").Value = 36000
$ws.Range("C14").Value = 1333
$ws.Range("D14").Value = 2667
$ws.Range("E14").Value = 24.08722222222222
$ws.Range("F14").Value = 25.09902475618905
$ws.Range("G14").Value = 23.66891638545182
$ws.Range("H14").Value = 1491
$ws.Range("I14").Value = 0.1719445533593191
$ws.Range("J14").Value = 972
$ws.Range("K14").Value = 2.905221627761007
$ws.Range("L14").Value = 950
$ws.Range("M14").Value = 1.504950495049505
$ws.Range("N14").Value = 7
$ws.Range("O14").Value = 0.0008072514242221556

# Row 16
$ws.Range("B16").Value = 4486
$ws.Range("C16").Value = 105
$ws.Range("D16").Value = 211
$ws.Range("E16").Value = 42.65202853321445
$ws.Range("F16").Value = 18.85714285714286
$ws.Range("G16").Value = 18.40284360189574
$ws.Range("H16").Value = 2780
$ws.Range("I16").Value = 1.452933828794222
$ws.Range("J16").Value = 651
$ws.Range("K16").Value = 32.87878787878788
$ws.Range("L16").Value = 911
$ws.Range("M16").Value = 23.4612413082668
$ws.Range("N16").Value = 2
$ws.Range("O16").Value = 0.001045276135823181

# Row 17
$ws.Range("B17").Value = 2468
$ws.Range("C17").Value = 353
$ws.Range("D17").Value = 706
$ws.Range("E17").Value = 24.32698541329011
$ws.Range("F17").Value = 22.99150141643059
$ws.Range("G17").Value = 26.16005665722379
$ws.Range("H17").Value = 1354
$ws.Range("I17").Value = 2.255200786155665
$ws.Range("J17").Value = 1212
$ws.Range("K17").Value = 14.93346476096599
$ws.Range("L17").Value = 1482
$ws.Range("M17").Value = 8.024256862851264

# Row 20
$ws.Range("B20").Value = 595
$ws.Range("C20").Value = 85
$ws.Range("D20").Value = 171
$ws.Range("E20").Value = 53.43697478991596
$ws.Range("F20").Value = 52.04705882352941
$ws.Range("G20").Value = 56.98245614035088
$ws.Range("H20").Value = 1527
$ws.Range("I20").Value = 4.802641924830948
$ws.Range("J20").Value = 703
$ws.Range("K20").Value = 15.89059674502712
$ws.Range("L20").Value = 981
$ws.Range("M20").Value = 10.06773399014778
$ws.Range("N20").Value = 10
$ws.Range("O20").Value = 0.03145148608271741
$ws.Range("P20").Value = 3
$ws.Range("Q20").Value = 0.06781193490054249
$ws.Range("R20").Value = 4
$ws.Range("S20").Value = 0.04105090311986864
